$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update open/close/high/low price, shares_outstanding, and fixed_ticker
# for each data row (rows 2-29, 31-40; row 30 already matches target state).
$ws.Cells.Item(2, 4).Value2 = 22.56767573870786
$ws.Cells.Item(2, 5).Value2 = 18.47954368591309
$ws.Cells.Item(2, 6).Value2 = 24.10222265071334
$ws.Cells.Item(2, 7).Value2 = 17.56529917706477
$ws.Cells.Item(2, 8).Value2 = 178998669
$ws.Cells.Item(2, 9).Value2 = "DXC"

$ws.Cells.Item(3, 4).Value2 = 25.75853278637041
$ws.Cells.Item(3, 5).Value2 = 23.263427734375
$ws.Cells.Item(3, 6).Value2 = 25.96040201767212
$ws.Cells.Item(3, 7).Value2 = 19.5974797165072
$ws.Cells.Item(3, 8).Value2 = 178998669
$ws.Cells.Item(3, 9).Value2 = "DXC"

$ws.Cells.Item(4, 4).Value2 = 40.30586368069592
$ws.Cells.Item(4, 5).Value2 = 38.88320922851562
$ws.Cells.Item(4, 6).Value2 = 40.79362805350653
$ws.Cells.Item(4, 7).Value2 = 38.37918006840941
$ws.Cells.Item(4, 8).Value2 = 178998669
$ws.Cells.Item(4, 9).Value2 = "DXC"

$ws.Cells.Item(5, 4).Value2 = 42.344859100305
$ws.Cells.Item(5, 5).Value2 = 44.39972305297852
$ws.Cells.Item(5, 6).Value2 = 46.4301209149399
$ws.Cells.Item(5, 7).Value2 = 41.10541643784433
$ws.Cells.Item(5, 8).Value2 = 178998669
$ws.Cells.Item(5, 9).Value2 = "DXC"

$ws.Cells.Item(6, 4).Value2 = 48.62341580256726
$ws.Cells.Item(6, 5).Value2 = 50.8384017944336
$ws.Cells.Item(6, 6).Value2 = 51.32063021952053
$ws.Cells.Item(6, 7).Value2 = 46.63728764391581
$ws.Cells.Item(6, 8).Value2 = 178998669
$ws.Cells.Item(6, 9).Value2 = "DXC"

$ws.Cells.Item(7, 4).Value2 = 56.77500184934564
$ws.Cells.Item(7, 5).Value2 = 61.70555114746094
$ws.Cells.Item(7, 6).Value2 = 63.87597595435853
$ws.Cells.Item(7, 7).Value2 = 55.49731980194213
$ws.Cells.Item(7, 8).Value2 = 178998669
$ws.Cells.Item(7, 9).Value2 = "DXC"

$ws.Cells.Item(8, 4).Value2 = 63.01753633352818
$ws.Cells.Item(8, 5).Value2 = 64.34750366210938
$ws.Cells.Item(8, 6).Value2 = 66.27678244066387
$ws.Cells.Item(8, 7).Value2 = 62.09805292711641
$ws.Cells.Item(8, 8).Value2 = 178998669
$ws.Cells.Item(8, 9).Value2 = "DXC"

$ws.Cells.Item(9, 4).Value2 = 70.65438563459843
$ws.Cells.Item(9, 5).Value2 = 75.29447174072266
$ws.Cells.Item(9, 6).Value2 = 76.84116710943074
$ws.Cells.Item(9, 7).Value2 = 69.8069908581845
$ws.Cells.Item(9, 8).Value2 = 178998669
$ws.Cells.Item(9, 9).Value2 = "DXC"

$ws.Cells.Item(10, 4).Value2 = 78.22628234508909
$ws.Cells.Item(10, 5).Value2 = 82.05928039550781
$ws.Cells.Item(10, 6).Value2 = 84.86190872334022
$ws.Cells.Item(10, 7).Value2 = 77.93777242516559
$ws.Cells.Item(10, 8).Value2 = 178998669
$ws.Cells.Item(10, 9).Value2 = "DXC"

$ws.Cells.Item(11, 4).Value2 = 83.101684696696
$ws.Cells.Item(11, 5).Value2 = 85.09996032714844
$ws.Cells.Item(11, 6).Value2 = 86.69362137342395
$ws.Cells.Item(11, 7).Value2 = 81.20250355402943
$ws.Cells.Item(11, 8).Value2 = 178998669
$ws.Cells.Item(11, 9).Value2 = "DXC"

$ws.Cells.Item(12, 4).Value2 = 76.22051588744115
$ws.Cells.Item(12, 5).Value2 = 81.07099914550781
$ws.Cells.Item(12, 6).Value2 = 84.69690389804451
$ws.Cells.Item(12, 7).Value2 = 76.22051588744115
$ws.Cells.Item(12, 8).Value2 = 178998669
$ws.Cells.Item(12, 9).Value2 = "DXC"

$ws.Cells.Item(13, 4).Value2 = 90.21391227710632
$ws.Cells.Item(13, 5).Value2 = 69.82231140136719
$ws.Cells.Item(13, 6).Value2 = 90.74120162745702
$ws.Cells.Item(13, 7).Value2 = 65.75740715251168
$ws.Cells.Item(13, 8).Value2 = 178998669
$ws.Cells.Item(13, 9).Value2 = "DXC"

$ws.Cells.Item(14, 4).Value2 = 50.27139579897369
$ws.Cells.Item(14, 5).Value2 = 61.6565055847168
$ws.Cells.Item(14, 6).Value2 = 61.94497515181015
$ws.Cells.Item(14, 7).Value2 = 50.11754316230884
$ws.Cells.Item(14, 8).Value2 = 178998669
$ws.Cells.Item(14, 9).Value2 = "DXC"

$ws.Cells.Item(15, 4).Value2 = 62.24394650558868
$ws.Cells.Item(15, 5).Value2 = 63.40125274658203
$ws.Cells.Item(15, 6).Value2 = 64.30781396207679
$ws.Cells.Item(15, 7).Value2 = 61.23130170522609
$ws.Cells.Item(15, 8).Value2 = 178998669
$ws.Cells.Item(15, 9).Value2 = "DXC"

$ws.Cells.Item(16, 4).Value2 = 54.26171959650711
$ws.Cells.Item(16, 5).Value2 = 54.02921295166016
$ws.Cells.Item(16, 6).Value2 = 55.44364007729632
$ws.Cells.Item(16, 7).Value2 = 52.68260057207323
$ws.Cells.Item(16, 8).Value2 = 178998669
$ws.Cells.Item(16, 9).Value2 = "DXC"

$ws.Cells.Item(17, 4).Value2 = 28.93650111172137
$ws.Cells.Item(17, 5).Value2 = 26.97685241699219
$ws.Cells.Item(17, 6).Value2 = 28.99499756034309
$ws.Cells.Item(17, 7).Value2 = 25.36818613313049
$ws.Cells.Item(17, 8).Value2 = 178998669
$ws.Cells.Item(17, 9).Value2 = "DXC"

$ws.Cells.Item(18, 4).Value2 = 37.29325533817211
$ws.Cells.Item(18, 5).Value2 = 31.26239776611328
$ws.Cells.Item(18, 6).Value2 = 37.62666877263113
$ws.Cells.Item(18, 7).Value2 = 30.77208411195568
$ws.Cells.Item(18, 8).Value2 = 178998669
$ws.Cells.Item(18, 9).Value2 = "DXC"

$ws.Cells.Item(19, 4).Value2 = 12.03999996185303
$ws.Cells.Item(19, 5).Value2 = 18.1299991607666
$ws.Cells.Item(19, 6).Value2 = 19.32999992370605
$ws.Cells.Item(19, 7).Value2 = 11.35000038146973
$ws.Cells.Item(19, 8).Value2 = 178998669
$ws.Cells.Item(19, 9).Value2 = "DXC"

$ws.Cells.Item(20, 4).Value2 = 16.45000076293945
$ws.Cells.Item(20, 5).Value2 = 17.90999984741211
$ws.Cells.Item(20, 6).Value2 = 17.98999977111816
$ws.Cells.Item(20, 7).Value2 = 14.72999954223633
$ws.Cells.Item(20, 8).Value2 = 178998669
$ws.Cells.Item(20, 9).Value2 = "DXC"

$ws.Cells.Item(21, 4).Value2 = 18.07999992370605
$ws.Cells.Item(21, 5).Value2 = 18.42000007629395
$ws.Cells.Item(21, 6).Value2 = 21.02000045776367
$ws.Cells.Item(21, 7).Value2 = 16.8799991607666
$ws.Cells.Item(21, 8).Value2 = 178998669
$ws.Cells.Item(21, 9).Value2 = "DXC"

$ws.Cells.Item(22, 4).Value2 = 26.95000076293945
$ws.Cells.Item(22, 5).Value2 = 28.20000076293945
$ws.Cells.Item(22, 6).Value2 = 30.13999938964844
$ws.Cells.Item(22, 7).Value2 = 25.55999946594238
$ws.Cells.Item(22, 8).Value2 = 178998669
$ws.Cells.Item(22, 9).Value2 = "DXC"

$ws.Cells.Item(23, 4).Value2 = 31.38999938964844
$ws.Cells.Item(23, 5).Value2 = 32.90999984741211
$ws.Cells.Item(23, 6).Value2 = 33.63000106811523
$ws.Cells.Item(23, 7).Value2 = 30.55999946594238
$ws.Cells.Item(23, 8).Value2 = 178998669
$ws.Cells.Item(23, 9).Value2 = "DXC"

$ws.Cells.Item(24, 4).Value2 = 39.22000122070312
$ws.Cells.Item(24, 5).Value2 = 39.97999954223633
$ws.Cells.Item(24, 6).Value2 = 40.97000122070312
$ws.Cells.Item(24, 7).Value2 = 36.75
$ws.Cells.Item(24, 8).Value2 = 178998669
$ws.Cells.Item(24, 9).Value2 = "DXC"

$ws.Cells.Item(25, 4).Value2 = 34
$ws.Cells.Item(25, 5).Value2 = 32.56999969482422
$ws.Cells.Item(25, 6).Value2 = 35.59000015258789
$ws.Cells.Item(25, 7).Value2 = 32.11000061035156
$ws.Cells.Item(25, 8).Value2 = 178998669
$ws.Cells.Item(25, 9).Value2 = "DXC"

$ws.Cells.Item(26, 4).Value2 = 32.40000152587891
$ws.Cells.Item(26, 5).Value2 = 30.07999992370605
$ws.Cells.Item(26, 6).Value2 = 34.63999938964844
$ws.Cells.Item(26, 7).Value2 = 28.56999969482422
$ws.Cells.Item(26, 8).Value2 = 178998669
$ws.Cells.Item(26, 9).Value2 = "DXC"

$ws.Cells.Item(27, 4).Value2 = 32.61999893188477
$ws.Cells.Item(27, 5).Value2 = 28.70000076293945
$ws.Cells.Item(27, 6).Value2 = 32.81000137329102
$ws.Cells.Item(27, 7).Value2 = 28.54000091552734
$ws.Cells.Item(27, 8).Value2 = 178998669
$ws.Cells.Item(27, 9).Value2 = "DXC"

$ws.Cells.Item(28, 4).Value2 = 30.17000007629395
$ws.Cells.Item(28, 5).Value2 = 31.60000038146973
$ws.Cells.Item(28, 6).Value2 = 31.98999977111816
$ws.Cells.Item(28, 7).Value2 = 27.31999969482422
$ws.Cells.Item(28, 8).Value2 = 178998669
$ws.Cells.Item(28, 9).Value2 = "DXC"

$ws.Cells.Item(29, 4).Value2 = 24.96999931335449
$ws.Cells.Item(29, 5).Value2 = 28.75
$ws.Cells.Item(29, 6).Value2 = 29.97999954223633
$ws.Cells.Item(29, 7).Value2 = 24.57999992370605
$ws.Cells.Item(29, 8).Value2 = 178998669
$ws.Cells.Item(29, 9).Value2 = "DXC"

$ws.Cells.Item(31, 4).Value2 = 25.6299991607666
$ws.Cells.Item(31, 5).Value2 = 23.85000038146973
$ws.Cells.Item(31, 6).Value2 = 25.86000061035156
$ws.Cells.Item(31, 7).Value2 = 22.73999977111816
$ws.Cells.Item(31, 8).Value2 = 178998669
$ws.Cells.Item(31, 9).Value2 = "DXC"

$ws.Cells.Item(32, 4).Value2 = 26.60000038146973
$ws.Cells.Item(32, 5).Value2 = 27.64999961853028
$ws.Cells.Item(32, 6).Value2 = 28.88999938964844
$ws.Cells.Item(32, 7).Value2 = 25.8700008392334
$ws.Cells.Item(32, 8).Value2 = 178998669
$ws.Cells.Item(32, 9).Value2 = "DXC"

$ws.Cells.Item(33, 4).Value2 = 20.77000045776367
$ws.Cells.Item(33, 5).Value2 = 20.17000007629395
$ws.Cells.Item(33, 6).Value2 = 22.63999938964844
$ws.Cells.Item(33, 7).Value2 = 19.29999923706055
$ws.Cells.Item(33, 8).Value2 = 178998669
$ws.Cells.Item(33, 9).Value2 = "DXC"

$ws.Cells.Item(34, 4).Value2 = 22.70000076293945
$ws.Cells.Item(34, 5).Value2 = 21.79999923706055
$ws.Cells.Item(34, 6).Value2 = 24.47999954223633
$ws.Cells.Item(34, 7).Value2 = 21.76000022888184
$ws.Cells.Item(34, 8).Value2 = 178998669
$ws.Cells.Item(34, 9).Value2 = "DXC"

$ws.Cells.Item(35, 4).Value2 = 21.19000053405762
$ws.Cells.Item(35, 5).Value2 = 19.48999977111816
$ws.Cells.Item(35, 6).Value2 = 22.14999961853028
$ws.Cells.Item(35, 7).Value2 = 19.47999954223633
$ws.Cells.Item(35, 8).Value2 = 178998669
$ws.Cells.Item(35, 9).Value2 = "DXC"

$ws.Cells.Item(36, 4).Value2 = 19.09000015258789
$ws.Cells.Item(36, 5).Value2 = 20.34000015258789
$ws.Cells.Item(36, 6).Value2 = 20.75
$ws.Cells.Item(36, 7).Value2 = 17.94000053405762
$ws.Cells.Item(36, 8).Value2 = 178998669
$ws.Cells.Item(36, 9).Value2 = "DXC"

$ws.Cells.Item(37, 4).Value2 = 20.64999961853028
$ws.Cells.Item(37, 5).Value2 = 19.86000061035156
$ws.Cells.Item(37, 6).Value2 = 21.35000038146973
$ws.Cells.Item(37, 7).Value2 = 19.65999984741211
$ws.Cells.Item(37, 8).Value2 = 178998669
$ws.Cells.Item(37, 9).Value2 = "DXC"

$ws.Cells.Item(38, 4).Value2 = 20.21999931335449
$ws.Cells.Item(38, 5).Value2 = 21.71999931335449
$ws.Cells.Item(38, 6).Value2 = 21.94000053405762
$ws.Cells.Item(38, 7).Value2 = 19.13999938964844
$ws.Cells.Item(38, 8).Value2 = 178998669
$ws.Cells.Item(38, 9).Value2 = "DXC"

$ws.Cells.Item(39, 4).Value2 = 16.90999984741211
$ws.Cells.Item(39, 5).Value2 = 15.52000045776367
$ws.Cells.Item(39, 6).Value2 = 17.68000030517578
$ws.Cells.Item(39, 7).Value2 = 13.4399995803833
$ws.Cells.Item(39, 8).Value2 = 178998669
$ws.Cells.Item(39, 9).Value2 = "DXC"

$ws.Cells.Item(40, 4).Value2 = 15.14000034332275
$ws.Cells.Item(40, 5).Value2 = 13.60999965667725
$ws.Cells.Item(40, 6).Value2 = 16.45000076293945
$ws.Cells.Item(40, 7).Value2 = 13.47999954223633
$ws.Cells.Item(40, 8).Value2 = 178998669
$ws.Cells.Item(40, 9).Value2 = "DXC"
